# Insert a new blank row above the old row 3 ("111"), pushing it and the
# row below ("1111") down by one. This grows the used range from A1:D4 to
# A1:D5 and leaves a blank row 3 in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(3).Insert()

# Move/leave the active selection on C4, matching the author's final
# cursor position after the edit.
$ws.Range("C4").Select()
